$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 3.05
$ws.Range("N3").Value = 1.65
$ws.Range("P3").Value = 1.65
$ws.Range("T3").Value = 1.01
$ws.Range("V3").Value = 1.45
$ws.Range("Z3").Value = 26
$ws.Range("K4").Value = 970
$ws.Range("N4").Value = 1.98
$ws.Range("P4").Value = 1.98
$ws.Range("U4").Value = 2.52
$ws.Range("N5").Value = 7.2
$ws.Range("O5").Value = 1.13
$ws.Range("R5").Value = 1.83
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.98
$ws.Range("U5").Value = 1.84
$ws.Range("X5").Value = 980
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AC5").Value = 980
$ws.Range("AD5").Value = 1000
$ws.Range("AH5").Value = 980
$ws.Range("AI5").Value = 980
$ws.Range("AK5").Value = 280
$ws.Range("AM5").Value = 190
$ws.Range("N6").Value = 1.02
$ws.Range("R6").Value = 1.81
$ws.Range("AL6").Value = 970
$ws.Range("AN6").Value = 5.4
$ws.Range("G7").Value = 3.6
$ws.Range("L7").Value = 1.29
$ws.Range("N7").Value = 5.9
$ws.Range("P7").Value = 2.62
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 2.48
$ws.Range("U7").Value = 2.68
$ws.Range("AO7").Value = 10
$ws.Range("F8").Value = 3.3
$ws.Range("H8").Value = 2.16
$ws.Range("I8").Value = 2.18
$ws.Range("L8").Value = 1.29
$ws.Range("N8").Value = 6
$ws.Range("P8").Value = 2.64
$ws.Range("U8").Value = 2.72
$ws.Range("V8").Value = 1.84
$ws.Range("X8").Value = 23
$ws.Range("AF8").Value = 27
$ws.Range("AK8").Value = 32
$ws.Range("AO8").Value = 10.5
$ws.Range("M9").Value = 1.17
$ws.Range("R9").Value = 1.13
$ws.Range("T9").Value = 2.62
$ws.Range("U9").Value = 1.53
$ws.Range("V9").Value = 1.19
$ws.Range("W9").Value = 1.94
$ws.Range("X9").Value = 6.6
$ws.Range("Y9").Value = 14.5
$ws.Range("Z9").Value = 46
$ws.Range("AA9").Value = 280
$ws.Range("AB9").Value = 5.5
$ws.Range("AC9").Value = 9.199999999999999
$ws.Range("AD9").Value = 28
$ws.Range("AE9").Value = 190
$ws.Range("AF9").Value = 10
$ws.Range("AG9").Value = 13
$ws.Range("AH9").Value = 970
$ws.Range("AI9").Value = 240
$ws.Range("AJ9").Value = 26
$ws.Range("AK9").Value = 36
$ws.Range("AL9").Value = 120
$ws.Range("AM9").Value = 490
$ws.Range("AN9").Value = 970
$ws.Range("AO9").Value = 420
$ws.Range("F10").Value = 2.54
$ws.Range("G10").Value = 2.92
$ws.Range("H10").Value = 2.78
$ws.Range("I10").Value = 3.25
$ws.Range("J10").Value = 3.1
$ws.Range("K10").Value = 3.65
$ws.Range("V10").Value = 1.45
$ws.Range("W10").Value = 1.52
$ws.Range("G11").Value = 3.8
$ws.Range("H11").Value = 2.06
$ws.Range("J11").Value = 3.7
$ws.Range("K11").Value = 4.1
$ws.Range("Q11").Value = 1.7
$ws.Range("T11").Value = 1.53
$ws.Range("W11").Value = 1.35
$ws.Range("AI11").Value = 970
$ws.Range("Q12").Value = 1.74
$ws.Range("R12").Value = 1.15
$ws.Range("S12").Value = 1.74
$ws.Range("U12").Value = 1.98
$ws.Range("Y12").Value = 970
$ws.Range("AB12").Value = 970
$ws.Range("AC12").Value = 970
$ws.Range("AD12").Value = 970
$ws.Range("AF12").Value = 970
$ws.Range("AG12").Value = 970
$ws.Range("AH12").Value = 970
$ws.Range("AJ12").Value = 970
$ws.Range("AK12").Value = 970
$ws.Range("AN12").Value = 970
$ws.Range("G13").Value = 2.76
$ws.Range("S13").Value = 1.74
$ws.Range("T13").Value = 1.62
$ws.Range("U13").Value = 2.2
$ws.Range("W13").Value = 1.57
$ws.Range("R14").Value = 1.47
$ws.Range("S14").Value = 2.26
$ws.Range("T14").Value = 1.64
$ws.Range("U14").Value = 1.87
$ws.Range("AN14").Value = 7.8
$ws.Range("F15").Value = 1.91
$ws.Range("G15").Value = 1.93
$ws.Range("I15").Value = 4.8
$ws.Range("J15").Value = 3.75
$ws.Range("K15").Value = 3.8
$ws.Range("L15").Value = 1.4
$ws.Range("P15").Value = 1.93
$ws.Range("Q15").Value = 2.02
$ws.Range("U15").Value = 2.04
$ws.Range("V15").Value = 1.26
$ws.Range("W15").Value = 2.08
$ws.Range("Y15").Value = 16
$ws.Range("AD15").Value = 18.5
$ws.Range("AG15").Value = 10
$ws.Range("AJ15").Value = 21
$ws.Range("AL15").Value = 38
$ws.Range("AN15").Value = 13.5
$ws.Range("G16").Value = 9.4
$ws.Range("I16").Value = 1.36
$ws.Range("J16").Value = 6.6
$ws.Range("K16").Value = 6.8
$ws.Range("L16").Value = 1.2
$ws.Range("U16").Value = 2.38
$ws.Range("V16").Value = 3.8
$ws.Range("W16").Value = 1.12
$ws.Range("AK16").Value = 110
$ws.Range("AM16").Value = 85
$ws.Range("AN16").Value = 70
$ws.Range("H17").Value = 2.06
$ws.Range("I17").Value = 2.08
$ws.Range("J17").Value = 3.9
$ws.Range("K17").Value = 3.95
$ws.Range("L17").Value = 1.28
$ws.Range("V17").Value = 1.93
$ws.Range("W17").Value = 1.34
$ws.Range("Y17").Value = 12
$ws.Range("AI17").Value = 28
$ws.Range("AJ17").Value = 70
$ws.Range("AL17").Value = 40
$ws.Range("J18").Value = 3.65
$ws.Range("L18").Value = 1.41
$ws.Range("V18").Value = 1.27
$ws.Range("W18").Value = 2.02
$ws.Range("G19").Value = 1.5
$ws.Range("L19").Value = 1.22
$ws.Range("N19").Value = 7.6
$ws.Range("S19").Value = 2.02
$ws.Range("V19").Value = 1.17
$ws.Range("W19").Value = 2.94
$ws.Range("F20").Value = 1.17
$ws.Range("G20").Value = 1.18
$ws.Range("H20").Value = 22
$ws.Range("J20").Value = 9.800000000000001
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = 1.19
$ws.Range("P20").Value = 3.35
$ws.Range("T20").Value = 2.18
$ws.Range("V20").Value = 1.04
$ws.Range("W20").Value = 6.8
$ws.Range("Y20").Value = 990
$ws.Range("AD20").Value = 80
$ws.Range("K21").Value = 13.5
$ws.Range("L21").Value = 1.14
$ws.Range("P21").Value = 4.5
$ws.Range("V21").Value = 1.03
$ws.Range("W21").Value = 8
$ws.Range("X21").Value = 1000
$ws.Range("AC21").Value = 36
$ws.Range("AN21").Value = 2.28
